$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are numeric-looking strings (e.g. "1.008"); Excel would
# otherwise auto-convert them to numbers, so we force Text format, write the
# value, then restore the cells original style to avoid leaving formatting
# changes behind.
$dRefs = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
$dVals = @("28.220.53", "1.881.79", "1.008", "314.43", "1.007", "0.5145", "0.3917", "0.08383", "41.61", "6.233", "20.73", "1.879.14", "7.266", "1.007", "0.00001103", "91.25", "0.06689", "17.81", "6.047", "28.268.21", "11.17", "2.272", "2.097.28", "159.90", "20.68", "125.33", "0.1061", "1.039", "5.872", "3.614", "9.684", "0.02449", "0.06577", "0.2190", "1.203", "0.6518", "5.006", "1.222", "11.35", "0.6149", "13.06", "1.287", "3.682", "1.232", "121.18", "0.06920")
for ($i = 0; $i -lt $dRefs.Length; $i++) {
    $cell = $ws.Range($dRefs[$i])
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $dVals[$i]
    $cell.Style = $origStyle
}

# Columns B, C, E are plain text (coin names, links, percentages) and do not
# risk being reinterpreted as numbers, so a direct Value assignment is safe.
$otherRefs = @("E2", "E3", "E4", "E5", "E6", "E7", "E8", "E9", "E10", "E11", "E12", "B13", "C13", "E13", "B14", "C14", "E14", "E15", "E16", "E17", "E18", "E19", "E20", "E21", "E22", "E23", "E24", "E25", "E26", "E27", "E28", "E29", "E30", "E31", "E32", "E33", "E34", "E35", "E36", "E37", "E38", "E39", "E40", "E41", "E42", "E43", "E44", "E45", "E46", "E47", "E48", "E49", "E50", "E51")
$otherVals = @("  +1.02%  ", "  +1.61%  ", "  +0.25%  ", "  +0.97%  ", "  +0.34%  ", "  +1.76%  ", "  +3.03%  ", "  +2.11%  ", "  +1.65%  ", "  +0.21%  ", "  +1.00%  ", "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "  +1.64%  ", "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "  +0.79%  ", "  +1.31%  ", "  +0.22%  ", "  +1.10%  ", "  +1.16%  ", "  +1.36%  ", "  +0.91%  ", "  +0.47%  ", "  +1.01%  ", "  +1.09%  ", "  +1.61%  ", "  +1.37%  ", "  +1.18%  ", "  -1.63%  ", "  +1.48%  ", "  +1.85%  ", "  +1.16%  ", "  +0.75%  ", "  +0.94%  ", "  +5.32%  ", "  +0.56%  ", "  +1.57%  ", "  +2.10%  ", "  +1.32%  ", "  +1.44%  ", "  +0.67%  ", "  +2.62%  ", "  +3.60%  ", "  -1.01%  ", "  +2.12%  ", "  +2.09%  ", "  -0.27%  ", "  +0.50%  ", "  +0.92%  ", "  +2.70%  ", "  +2.59%  ", "  +0.65%  ", "  +1.20%  ")
for ($i = 0; $i -lt $otherRefs.Length; $i++) {
    $ws.Range($otherRefs[$i]).Value = $otherVals[$i]
}